$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.033.68"
$ws.Range("E2").Value = "  -0.48%  "
$ws.Range("D3").Value = "2.550.99"
$ws.Range("E3").Value = "  +0.14%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'584.57"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.31%  "
$ws.Range("D6").Value = "'147.19"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.50%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("E8").Value = "  -1.02%  "
$ws.Range("E9").Value = "  -0.43%  "
$ws.Range("E10").Value = "  -3.33%  "
$ws.Range("E11").Value = "  -0.13%  "
$ws.Range("E12").Value = "  -1.05%  "
$ws.Range("E13").Value = "  -4.25%  "
$ws.Range("D14").Value = "3.007.89"
$ws.Range("E14").Value = "  +0.14%  "
$ws.Range("D15").Value = "62.945.34"
$ws.Range("E15").Value = "  -0.51%  "
$ws.Range("E16").Value = "  -0.71%  "
$ws.Range("D17").Value = "2.548.17"
$ws.Range("E17").Value = "  -0.02%  "
$ws.Range("D18").Value = "'11.36"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.79%  "
$ws.Range("D19").Value = "'336.38"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.30%  "
$ws.Range("E20").Value = "  -1.06%  "
$ws.Range("E21").Value = "  -1.68%  "
$ws.Range("E22").Value = "  -0.08%  "
$ws.Range("D23").Value = "'65.90"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.14%  "
$ws.Range("E24").Value = "  -0.47%  "
$ws.Range("E25").Value = "  +0.34%  "
$ws.Range("E26").Value = "  +0.13%  "
$ws.Range("D27").Value = "'1.48"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.40%  "
$ws.Range("D28").Value = "'8.38"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.38%  "
$ws.Range("D29").Value = "'7.42"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.97%  "
$ws.Range("E30").Value = "  +3.12%  "
$ws.Range("D31").Value = "0.0₃0814"
$ws.Range("E31").Value = "  -2.73%  "
$ws.Range("D32").Value = "'177.71"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.60%  "
$ws.Range("E33").Value = "  -1.08%  "
$ws.Range("D34").Value = "'416.14"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.19%  "
$ws.Range("D35").Value = "'19.15"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.01%  "
$ws.Range("D36").Value = "'0.400"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.95%  "
$ws.Range("E37").Value = "  +0.02%  "
$ws.Range("D38").Value = "'4.35"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.39%  "
$ws.Range("E39").Value = "  -1.14%  "
$ws.Range("E40").Value = "  +0.01%  "
$ws.Range("D41").Value = "'39.59"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.00%  "
$ws.Range("D42").Value = "'151.07"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.18%  "
$ws.Range("D43").Value = "'3.76"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.29%  "
$ws.Range("D44").Value = "'20.87"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.35%  "
$ws.Range("D45").Value = "'0.0540"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.42%  "
$ws.Range("D46").Value = "'0.601"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.31%  "
$ws.Range("D47").Value = "'0.0971"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.47%  "
$ws.Range("E48").Value = "  +0.04%  "
$ws.Range("D49").Value = "'18.31"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.86%  "
$ws.Range("D50").Value = "'1.71"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -5.94%  "
$ws.Range("E51").Value = "  +0.00%  "
